$wb = $excel.ActiveWorkbook

# Switch to the ETLE-capacity sheet (third sheet, 0-indexed tab 2)
$ws = $wb.Worksheets.Item("ETLE-capacity")
$ws.Activate()

# Update calibrated value in B2 from -0.3 to -90
$ws.Range("B2").Value = -90

# Update the active selection to B2 (previously B3)
$ws.Range("B2").Select()
